$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 114 / 115 - the two rows have their data (everything except the
#    running index in column A) swapped between them.
# ---------------------------------------------------------------------------

# New row 114 <- old row 115 data
$ws.Cells.Item(114, 2).Value  = 7559468
$ws.Cells.Item(114, 6).Value  = "Liverpool Montevideo"
$ws.Cells.Item(114, 7).Value  = "CA River Plate"
$ws.Cells.Item(114, 8).Value  = 2
$ws.Cells.Item(114, 9).Value  = 1
$ws.Cells.Item(114, 10).Value = "H"
$ws.Cells.Item(114, 11).Value = 1.7
$ws.Cells.Item(114, 12).Value = 3
$ws.Cells.Item(114, 13).Value = 5.75
$ws.Cells.Item(114, 14).Value = 1.833
$ws.Cells.Item(114, 15).Value = 3.2
$ws.Cells.Item(114, 16).Value = 4.5
$ws.Cells.Item(114, 17).Value = -0.5
$ws.Cells.Item(114, 18).Value = 1.925
$ws.Cells.Item(114, 19).Value = 1.925
$ws.Cells.Item(114, 20).Value = 2.25
$ws.Cells.Item(114, 21).Value = 2.025
$ws.Cells.Item(114, 22).Value = 1.825
$ws.Cells.Item(114, 23).Value = 0.833
$ws.Cells.Item(114, 24).Value = -1
$ws.Cells.Item(114, 25).Value = -1
$ws.Cells.Item(114, 26).Value = 0.925
$ws.Cells.Item(114, 27).Value = -1
$ws.Cells.Item(114, 28).Value = 1.025
$ws.Cells.Item(114, 29).Value = -1

# New row 115 <- old row 114 data
$ws.Cells.Item(115, 2).Value  = 7559469
$ws.Cells.Item(115, 6).Value  = "Montevideo Wanderers"
$ws.Cells.Item(115, 7).Value  = "Penarol"
$ws.Cells.Item(115, 8).Value  = 0
$ws.Cells.Item(115, 9).Value  = 0
$ws.Cells.Item(115, 10).Value = "D"
$ws.Cells.Item(115, 11).Value = 4.75
$ws.Cells.Item(115, 12).Value = 3.4
$ws.Cells.Item(115, 13).Value = 1.7
$ws.Cells.Item(115, 14).Value = 2.7
$ws.Cells.Item(115, 15).Value = 3.2
$ws.Cells.Item(115, 16).Value = 2.45
$ws.Cells.Item(115, 17).Value = 0
$ws.Cells.Item(115, 18).Value = 2.05
$ws.Cells.Item(115, 19).Value = 1.8
$ws.Cells.Item(115, 20).Value = 2.5
$ws.Cells.Item(115, 21).Value = 1.975
$ws.Cells.Item(115, 22).Value = 1.875
$ws.Cells.Item(115, 23).Value = -1
$ws.Cells.Item(115, 24).Value = 2.2
$ws.Cells.Item(115, 25).Value = -1
$ws.Cells.Item(115, 26).Value = 0
$ws.Cells.Item(115, 27).Value = 0
$ws.Cells.Item(115, 28).Value = -1
$ws.Cells.Item(115, 29).Value = 0.875

# ---------------------------------------------------------------------------
# 2) Row 177 - overwritten with a new fixture's odds (updated match).
# ---------------------------------------------------------------------------
$ws.Cells.Item(177, 2).Value  = 8051186
$ws.Cells.Item(177, 5).Value  = 45395.83333333334
$ws.Cells.Item(177, 6).Value  = "Deportivo Maldonado"
$ws.Cells.Item(177, 7).Value  = "Cerro Largo"
$ws.Cells.Item(177, 11).Value = 2.3
$ws.Cells.Item(177, 12).Value = 3
$ws.Cells.Item(177, 13).Value = 3.3
$ws.Cells.Item(177, 14).Value = 2.4
$ws.Cells.Item(177, 15).Value = 3
$ws.Cells.Item(177, 16).Value = 3.1
$ws.Cells.Item(177, 17).Value = -0.25
$ws.Cells.Item(177, 18).Value = 2.1
$ws.Cells.Item(177, 19).Value = 1.775
$ws.Cells.Item(177, 20).Value = 2.25
$ws.Cells.Item(177, 21).Value = 2.1
$ws.Cells.Item(177, 22).Value = 1.775

# ---------------------------------------------------------------------------
# 3) Row 178 (new) - the fixture that used to live in row 177 (same teams /
#    date) gets re-added with refreshed closing odds.
# ---------------------------------------------------------------------------

# Bring over the same cell formatting (style) used by every other data row
# for the running-index column (A, bold/bordered) and the date column
# (E, date number format) before filling in the values.
$ws.Range("A177").Copy()
$ws.Range("A178").PasteSpecial(-4122)
$ws.Range("A179").PasteSpecial(-4122)
$ws.Range("E177").Copy()
$ws.Range("E178").PasteSpecial(-4122)
$ws.Range("E179").PasteSpecial(-4122)

$ws.Cells.Item(178, 1).Value  = 176
$ws.Cells.Item(178, 2).Value  = 8051187
$ws.Cells.Item(178, 3).Value  = "Uruguay Primera División"
$ws.Cells.Item(178, 4).Value  = "Uruguay Apertura"
$ws.Cells.Item(178, 5).Value  = 45396.625
$ws.Cells.Item(178, 6).Value  = "Defensor Sporting"
$ws.Cells.Item(178, 7).Value  = "CA River Plate"
$ws.Cells.Item(178, 11).Value = 1.727
$ws.Cells.Item(178, 12).Value = 3.5
$ws.Cells.Item(178, 13).Value = 5
$ws.Cells.Item(178, 14).Value = 1.615
$ws.Cells.Item(178, 15).Value = 3.6
$ws.Cells.Item(178, 16).Value = 6
$ws.Cells.Item(178, 17).Value = -0.75
$ws.Cells.Item(178, 18).Value = 1.85
$ws.Cells.Item(178, 19).Value = 2
$ws.Cells.Item(178, 20).Value = 2.25
$ws.Cells.Item(178, 21).Value = 1.85
$ws.Cells.Item(178, 22).Value = 2
$ws.Cells.Item(178, 23).Value = 0
$ws.Cells.Item(178, 24).Value = 0
$ws.Cells.Item(178, 25).Value = 0
$ws.Cells.Item(178, 26).Value = 0
$ws.Cells.Item(178, 27).Value = 0

# ---------------------------------------------------------------------------
# 4) Row 179 (new) - a brand-new fixture.
# ---------------------------------------------------------------------------
$ws.Cells.Item(179, 1).Value  = 177
$ws.Cells.Item(179, 2).Value  = 8050912
$ws.Cells.Item(179, 3).Value  = "Uruguay Primera División"
$ws.Cells.Item(179, 4).Value  = "Uruguay Apertura"
$ws.Cells.Item(179, 5).Value  = 45397.75
$ws.Cells.Item(179, 6).Value  = "Montevideo Wanderers"
$ws.Cells.Item(179, 7).Value  = "Liverpool Montevideo"
$ws.Cells.Item(179, 11).Value = 3.2
$ws.Cells.Item(179, 12).Value = 3.3
$ws.Cells.Item(179, 13).Value = 2.2
$ws.Cells.Item(179, 14).Value = 3.4
$ws.Cells.Item(179, 15).Value = 3.3
$ws.Cells.Item(179, 16).Value = 2.1
$ws.Cells.Item(179, 17).Value = 0.25
$ws.Cells.Item(179, 18).Value = 2
$ws.Cells.Item(179, 19).Value = 1.85
$ws.Cells.Item(179, 20).Value = 2.25
$ws.Cells.Item(179, 21).Value = 1.9
$ws.Cells.Item(179, 22).Value = 1.95
$ws.Cells.Item(179, 23).Value = 0
$ws.Cells.Item(179, 24).Value = 0
$ws.Cells.Item(179, 25).Value = 0
$ws.Cells.Item(179, 26).Value = 0
$ws.Cells.Item(179, 27).Value = 0
